# Auto-generated edit script: updates cryptocurrency price/volume data
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells below hold plain-decimal price strings (e.g. "233.71") that must
# stay TEXT (matching the source data feed), so force Text number format
# before assigning -- otherwise COM auto-converts them to numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.477.79"
$ws.Range("E2").Value = "  -2.93%  "

$ws.Range("D3").Value = "2.247.06"
$ws.Range("E3").Value = "  -3.81%  "

$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "233.71"
$ws.Range("E5").Value = "  -1.99%  "

$ws.Range("D6").Value = "0.633"
$ws.Range("E6").Value = "  -4.11%  "

$ws.Range("D7").Value = "69.56"
$ws.Range("E7").Value = "  -2.92%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  -3.24%  "

$ws.Range("D10").Value = "0.0994"
$ws.Range("E10").Value = "  +1.31%  "

$ws.Range("D11").Value = "58.73"
$ws.Range("E11").Value = "  +1.11%  "

$ws.Range("D12").Value = "36.69"
$ws.Range("E12").Value = "  +13.64%  "

$ws.Range("E13").Value = "  -1.93%  "

$ws.Range("E14").Value = "  -5.07%  "

$ws.Range("D15").Value = "2.582.68"
$ws.Range("E15").Value = "  -3.71%  "

$ws.Range("D16").Value = "15.11"
$ws.Range("E16").Value = "  -5.53%  "

$ws.Range("D17").Value = "0.858"
$ws.Range("E17").Value = "  -3.56%  "

$ws.Range("D18").Value = "2.249.65"
$ws.Range("E18").Value = "  -3.79%  "

$ws.Range("D19").Value = "42.293.30"
$ws.Range("E19").Value = "  -2.99%  "

$ws.Range("D20").Value = "0.0₃0977"
$ws.Range("E20").Value = "  -2.50%  "

$ws.Range("E21").Value = "  -4.58%  "

$ws.Range("D22").Value = "73.48"
$ws.Range("E22").Value = "  -5.34%  "

$ws.Range("D23").Value = "236.41"
$ws.Range("E23").Value = "  -5.43%  "

$ws.Range("D24").Value = "1.99"
$ws.Range("E24").Value = "  +4.61%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("E26").Value = "  -0.95%  "

$ws.Range("E27").Value = "  -3.38%  "

$ws.Range("D28").Value = "10.01"
$ws.Range("E28").Value = "  -2.36%  "

$ws.Range("E29").Value = "  -2.13%  "

$ws.Range("D30").Value = "169.93"
$ws.Range("E30").Value = "  -3.00%  "

$ws.Range("D31").Value = "20.57"
$ws.Range("E31").Value = "  -6.61%  "

$ws.Range("E32").Value = "  -3.68%  "

$ws.Range("E33").Value = "  -5.16%  "

$ws.Range("D34").Value = "0.0732"
$ws.Range("E34").Value = "  +0.14%  "

$ws.Range("D35").Value = "5.34"
$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("D36").Value = "4.70"
$ws.Range("E36").Value = "  -6.44%  "

$ws.Range("D37").Value = "3.66"
$ws.Range("E37").Value = "  -1.74%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0281"
$ws.Range("E38").Value = "  +5.05%  "

$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "22.01"
$ws.Range("E39").Value = "  +17.87%  "

$ws.Range("E40").Value = "  -3.30%  "

$ws.Range("E41").Value = "  -5.32%  "

$ws.Range("D42").Value = "65.61"
$ws.Range("E42").Value = "  +1.34%  "

$ws.Range("D43").Value = "9.20"
$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("D44").Value = "4.92"
$ws.Range("E44").Value = "  -10.59%  "

$ws.Range("E45").Value = "  -1.72%  "

$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("D48").Value = "4.54"
$ws.Range("E48").Value = "  +11.95%  "

$ws.Range("B49").Value = "Celestia"
$ws.Range("C49").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D49").Value = "10.16"
$ws.Range("E49").Value = "  +9.57%  "

$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "1.18"
$ws.Range("E50").Value = "  -2.66%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").Value = "  -2.80%  "
